$wb = $excel.ActiveWorkbook

# Update values on the "展览" sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 140
$ws1.Range("F5").Value = 2984
$ws1.Range("F6").Value = 303

# Update the same values on the "全部类型" sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 140
$ws4.Range("F5").Value = 2984
$ws4.Range("F6").Value = 303
